# Update "想去人数" (column F) counts that changed between scrape runs.
# Sheet "展览" (Worksheets index/name) rows: F5, F12, F20, F25, F29, F40, F42
# Sheet "全部类型" rows: F5, F12, F20, F26, F30, F42, F44

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 7491
$ws1.Range("F12").Value = 4335
$ws1.Range("F20").Value = 512
$ws1.Range("F25").Value = 1698
$ws1.Range("F29").Value = 111
$ws1.Range("F40").Value = 92
$ws1.Range("F42").Value = 48

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 7491
$ws4.Range("F12").Value = 4335
$ws4.Range("F20").Value = 512
$ws4.Range("F26").Value = 1698
$ws4.Range("F30").Value = 111
$ws4.Range("F42").Value = 92
$ws4.Range("F44").Value = 48
